# REPORTE.xlsx update
# - Appends "(VOLARSE DE NUEVO)" to every cell that currently reads exactly
#   "FALTA LINEAS DE VUELO" or "FALTA POSTPROCESO DE LA BASE DEL VUELO",
#   EXCEPT for a couple of rows that instead get a more specific,
#   brand-new observation text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldLinea     = "FALTA LINEAS DE VUELO"
$newLinea     = "FALTA LINEAS DE VUELO (VOLARSE DE NUEVO)"

$oldPost      = "FALTA POSTPROCESO DE LA BASE DEL VUELO"
$newPost      = "FALTA POSTPROCESO DE LA BASE DEL VUELO (VOLARSE DE NUEVO)"

# Row D80 (I11D71b4a2a) gets a distinct, more specific observation instead
# of the generic "(VOLARSE DE NUEVO)" suffix.
$row80Text    = "FALTA LINEAS DE VUELO (CON RECORTES DE AREAS CONTIGUAS SE OBTIENEN AREAS FALTANTES)"

# Rows D84:D87 (I11D71b2d3c, I11D71b2d3d, I11D71b4b1a, I11D71b4b1b) get a
# distinct, more specific observation instead of the generic suffix.
$postSpecialText = "FALTA POSTPROCESO DE LA BASE DEL VUELO (ACCESO PENDIENTE)"

$postSpecialRows = @(84, 85, 86, 87)

$lastRow = 203
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2

    if ($r -eq 80) {
        if ($val -eq $oldLinea) {
            $cell.Value2 = $row80Text
        }
        continue
    }

    if ($postSpecialRows -contains $r) {
        if ($val -eq $oldPost) {
            $cell.Value2 = $postSpecialText
        }
        continue
    }

    if ($val -eq $oldLinea) {
        $cell.Value2 = $newLinea
    }
    elseif ($val -eq $oldPost) {
        $cell.Value2 = $newPost
    }
}

# Column D widened (bestFit) now that it holds longer text.
$ws.Columns.Item(4).AutoFit()
